# "add 1 scene quest" - Scene.xlsx
# Updates QuestDungeon ("H" column) / Quest ("F" column) tag lists on the
# "Scene" sheet's data table to add a new scene quest (suntemple;2), and
# redistribute a few quest tags between rows 23/25 and 28/29/30.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 (13010005 / 布萨特高塔): Quest (F) was empty, now gets a value.
$ws.Range("F7").Value = "suntemple;2"

# Row 23 (13020001): Quest (F) was empty, now gets a value; QuestDungeon (H)
# loses the "portal;1" tag (it moves to row 25's neighbours below).
$ws.Range("F23").Value = "colorpool;1|barn;1|portal;1"
$ws.Range("H23").Value = "trees;4|manflower;2|sandland;2|cliff;2"

# Row 25 (13020011): QuestDungeon (H) tag list removed entirely.
$ws.Range("H25").Value = $null

# Row 28 (13020021): Quest (F) gains "brokehouse;2"; QuestDungeon (H) loses
# the "brokehouse;3" tag, keeping only "trees;4".
$ws.Range("F28").Value = "sandflow;2|brokehouse;2"
$ws.Range("H28").Value = "trees;4"

# Row 29 (13020022): Quest (F) gains "brokehouse;1"; QuestDungeon (H) gets a
# brand-new tag.
$ws.Range("F29").Value = "barn;1|diarybook;1|brokehouse;1"
$ws.Range("H29").Value = "suntemple;2"

# Row 30 (13020023): QuestDungeon (H) gets a brand-new tag.
$ws.Range("H30").Value = "suntemple;2"

# Selection moved to F23 as the last-edited/highlighted cell.
$ws.Range("F23").Select()
